$d = $word.ActiveDocument

# Locate the "© 2020 ... Jekyll ..." footer paragraph via Find.
$rng = $d.Content
$found = $rng.Find.Execute("© 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the copyright footer paragraph."
}

# Figure out which paragraph (1-based index) contains the match.
$all = $d.Paragraphs
$copyIndex = -1
for ($i = 1; $i -le $all.Count; $i++) {
    $p = $all.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $copyIndex = $i
        break
    }
}
if ($copyIndex -lt 3) {
    throw "Unexpected document layout around the copyright paragraph."
}

# The two paragraphs immediately before the copyright line are the blank
# "Normal" paragraph and the blank "pageBreakBefore" paragraph that were
# inserted together with it. Remove all three (from the start of the first
# blank paragraph through the end of the copyright paragraph).
$startPara = $all.Item($copyIndex - 2)
$endPara = $all.Item($copyIndex)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
